$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...architettura client-server." -> "...architettura client-server. "
# (add a trailing space to the run that holds the full stop after the bold
#  "client-server" run; keep the bold run untouched)
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.ClearFormatting()
$found1 = $r1.Find.Execute("client-server")
if ($found1) {
    $dotRange = $d.Range($r1.End, $r1.End + 1)
    if ($dotRange.Text -eq ".") {
        $dotRange.Text = ". "
    }
}

# ---------------------------------------------------------------------------
# Change 2: rewording of the smart-contract recompilation sentence
# ---------------------------------------------------------------------------
$old2 = ", integrabile con il nostro server Node.js. Lo smart contract viene ricompilato ad ogni riavvio del server in maniera automatica."
$new2 = ", integrabile con il nostro server Node.js. Lo smart contract viene ricompilato ogni volta prima di fare il deploy del contratto. "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: shorten the "registro parallelo" paragraph, dropping the bold
# "off-chain" run and the trailing examples/ellipsis runs entirely.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.ClearFormatting()
$found3 = $r3.Find.Execute("Abbiamo deciso di mantenere anche un registro parallelo in cui andiamo a memorizzare le operazioni che avvengono ")
if ($found3) {
    $p3start = $r3.Start

    $r3b = $d.Range($r3.End, $d.Content.End)
    $r3b.Find.ClearFormatting()
    $found3b = $r3b.Find.Execute("…")
    if ($found3b) {
        $p3end = $r3b.End
        $fullRange = $d.Range($p3start, $p3end)
        $fullRange.Text = "Abbiamo deciso di mantenere anche un registro parallelo in cui andiamo a memorizzare tutti gli accessi al sito web."
    }
}

# ---------------------------------------------------------------------------
# Change 4: "autenticazione ed autenticazione" -> "autorizzazione ed autenticazione"
# (the whole sentence lives in a single run in the source document, so a
#  plain Find/Replace over that exact run text keeps everything tidy)
# ---------------------------------------------------------------------------
$old4 = "Per proteggere il server potrebbe essere utile inserire procedure di autenticazione ed autenticazione al momento dell"
$new4 = "Per proteggere il server potrebbe essere utile inserire procedure di autorizzazione ed autenticazione al momento dell"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null


Write-Output "All edits applied."
